$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.044.59"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.02%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.832.92"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.15%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.92%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6167"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.56%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.002"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.17%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07447"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.30%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2924"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.69%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.02"
$ws.Range("D10").Style = "Normal"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07681"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.36%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.828.61"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.14%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.998"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.09%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6714"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.07%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "82.56"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.77%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009143"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.06%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.890"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.20%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.041.87"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.13%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.081.73"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.13%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "232.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.25%  "

$ws.Range("E21").Value = "  +0.54%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.002"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.25%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.189"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.47%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.08%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.07"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.65%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1409"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.50%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.469"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.57%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.80"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.87%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.499"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.45%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.150"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.05%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.097"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.47%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.05521"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.54%  "

$ws.Range("E33").Value = "  +0.44%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.831"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.53%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7356"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.30%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.137"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.31%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.660"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.06%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.771"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.42%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01778"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.48%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.208.40"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.94%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.451"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.34%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8893"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.43%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.001"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.13%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.85"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.41%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.977.10"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.25%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "65.35"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.40%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000122"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.53%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5089"
$ws.Range("D48").Style = "Normal"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4065"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.14%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.113"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.24%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05817"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.48%  "
